$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Library")
$ws2 = $wb.Worksheets.Item("Library_Formula")

function Set-IndicatorRow($row, $name) {
    $ws2.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws2.Cells.Item($row, 1).Font.Name = "Trebuchet MS"
    $ws2.Cells.Item($row, 1).Font.Size = 10
    $ws2.Cells.Item($row, 2).Value = "LIB_EWS_RETAIL"
    $ws2.Cells.Item($row, 2).Font.Name = "Trebuchet MS"
    $ws2.Cells.Item($row, 2).Font.Size = 10
    $ws2.Cells.Item($row, 3).Value = $name
    $ws2.Cells.Item($row, 5).Value = "String"
    $ws2.Cells.Item($row, 5).Font.Name = "Trebuchet MS"
    $ws2.Cells.Item($row, 5).Font.Size = 10
    $ws2.Cells.Item($row, 6).Value = "String"
    $ws2.Cells.Item($row, 6).Font.Name = "Trebuchet MS"
    $ws2.Cells.Item($row, 6).Font.Size = 10
}

Set-IndicatorRow 25 "INDICATOR_52 "
Set-IndicatorRow 26 "INDICATOR_53 "
Set-IndicatorRow 27 "INDICATOR_54 "
Set-IndicatorRow 28 "INDICATOR_55 "
Set-IndicatorRow 29 "INDICATOR_60 "
Set-IndicatorRow 30 "INDICATOR_61 "
Set-IndicatorRow 31 "INDICATOR_62 "
Set-IndicatorRow 32 "INDICATOR_63 "
Set-IndicatorRow 33 "INDICATOR_64 "
Set-IndicatorRow 34 "INDICATOR_70 "
Set-IndicatorRow 35 "INDICATOR_71 "
Set-IndicatorRow 36 "INDICATOR_80 "
Set-IndicatorRow 37 "INDICATOR_81 "
Set-IndicatorRow 38 "INDICATOR_82 "
Set-IndicatorRow 39 "INDICATOR_83 "
Set-IndicatorRow 40 "INDICATOR_84 "
Set-IndicatorRow 41 "INDICATOR_85 "
Set-IndicatorRow 42 "INDICATOR_86 "
Set-IndicatorRow 43 "INDICATOR_87 "
Set-IndicatorRow 44 "INDICATOR_100"
Set-IndicatorRow 45 "INDICATOR_114"
Set-IndicatorRow 46 "INDICATOR_201"
Set-IndicatorRow 47 "INDICATOR_202"
Set-IndicatorRow 48 "INDICATOR_203"
Set-IndicatorRow 49 "INDICATOR_204"
Set-IndicatorRow 50 "INDICATOR_205"
Set-IndicatorRow 51 "INDICATOR_206"
Set-IndicatorRow 52 "INDICATOR_207"
Set-IndicatorRow 53 "INDICATOR_208"
Set-IndicatorRow 54 "INDICATOR_209"
Set-IndicatorRow 55 "INDICATOR_210"
Set-IndicatorRow 56 "INDICATOR_211"
Set-IndicatorRow 57 "INDICATOR_212"
Set-IndicatorRow 58 "INDICATOR_213"
Set-IndicatorRow 59 "INDICATOR_214"
Set-IndicatorRow 60 "INDICATOR_215"
Set-IndicatorRow 61 "INDICATOR_216"
Set-IndicatorRow 62 "INDICATOR_217"
Set-IndicatorRow 63 "INDICATOR_218"
Set-IndicatorRow 64 "INDICATOR_219"
Set-IndicatorRow 65 "INDICATOR_220"
Set-IndicatorRow 66 "INDICATOR_221"
Set-IndicatorRow 67 "INDICATOR_222"
Set-IndicatorRow 68 "INDICATOR_223"
Set-IndicatorRow 69 "INDICATOR_224"
Set-IndicatorRow 70 "INDICATOR_225"
Set-IndicatorRow 71 "INDICATOR_226"
Set-IndicatorRow 72 "INDICATOR_227"
Set-IndicatorRow 73 "INDICATOR_228"
Set-IndicatorRow 74 "INDICATOR_229"
Set-IndicatorRow 75 "INDICATOR_230"
Set-IndicatorRow 76 "INDICATOR_231"
Set-IndicatorRow 77 "INDICATOR_106"
Set-IndicatorRow 78 "INDICATOR_107"
Set-IndicatorRow 79 "INDICATOR_108"
Set-IndicatorRow 80 "INDICATOR_109"
Set-IndicatorRow 81 "INDICATOR_110"
Set-IndicatorRow 82 "INDICATOR_111"
Set-IndicatorRow 83 "INDICATOR_112"
Set-IndicatorRow 84 "INDICATOR_113"
Set-IndicatorRow 85 "INDICATOR_115"
Set-IndicatorRow 86 "INDICATOR_116"
Set-IndicatorRow 87 "INDICATOR_117"

# --- View state changes ---
# Library_Formula was the active sheet while its view was adjusted
# (zoom to 70%, selection moved to A10), then the user switched back
# to the Library sheet before saving, so Library ends up tabSelected.
$ws2.Activate()
$ws2.Range("A10").Select()
$excel.ActiveWindow.Zoom = 70

$ws1.Activate()
$ws1.Range("B2").Select()

Write-Host "Edit complete"
